$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Updated "last refreshed" note
$ws.Range("A1").Value = "Datos actualizados a 8 de Abril de 2020 a las 13:52"

# Swap the "Arroyo de la Luz" / "La Gomera" rows
$ws.Range("A62").Value = "Arroyo de la Luz"
$ws.Range("A63").Value = "La Gomera"

# Row 32 - Tenerife
$ws.Range("B32").Value = 1154
$ws.Range("C32").Value = 207
$ws.Range("D32").Value = 887
$ws.Range("E32").Value = 60

# Row 47 - Gran Canaria
$ws.Range("B47").Value = 446
$ws.Range("C47").Value = 43
$ws.Range("D47").Value = 377
$ws.Range("E47").Value = 26

# Row 56 - La Palma
$ws.Range("B56").Value = 68
$ws.Range("C56").Value = 8
$ws.Range("D56").Value = 57

# Row 57 - Lanzarote
$ws.Range("B57").Value = 63
$ws.Range("C57").Value = 2
$ws.Range("D57").Value = 59

# Row 59 - Fuerteventura
$ws.Range("C59").Value = 7
$ws.Range("D59").Value = 16

# Row 62 - now "Arroyo de la Luz"
$ws.Range("C62").Value = 0
$ws.Range("D62").Value = 7

# Row 63 - now "La Gomera"
$ws.Range("C63").Value = 2
$ws.Range("D63").Value = 5

# Row 64 - El Hierro
$ws.Range("C64").Value = 1
$ws.Range("D64").Value = 0
